$wb = $excel.ActiveWorkbook

$wsInputs = $wb.Worksheets.Item("Inputs")
$wsOptions = $wb.Worksheets.Item("Options")
$wsExport = $wb.Worksheets.Item("export_ready_inputs")

# Switch the "Payload/MTOW" selector from MTOW to Payload Weight
$wsInputs.Range("B7").Value = "Payload Weight"

# Update the associated value (payload weight target) from 20 to 4
$wsInputs.Range("C7").Value = 4

# Disable the hand-launchable requirement
$wsInputs.Range("C10").Value = $false

# Update selections / active sheet so the workbook reopens on Inputs
$wsInputs.Range("D11").Select()
$wsOptions.Range("B16").Select()

$wsInputs.Activate()
$wsInputs.Range("D11").Select()
